$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 14-17 (sending cluster "Resolving-Mac" block removed)
$ws.Rows("14:17").Delete()

# Update changed numeric values in rows 2-13
$ws.Range("G2").Value2 = 1.167885
$ws.Range("H2").Value2 = 3.503655
$ws.Range("I2").Value2 = 0.484739710372708
$ws.Range("J2").Value2 = 0.484739710372708
$ws.Range("M2").Value2 = 3.063353333333333
$ws.Range("N2").Value2 = 9.190059999999999
$ws.Range("O2").Value2 = 0.1884019917097105
$ws.Range("P2").Value2 = 0.1884019917097105
$ws.Range("Q2").Value2 = 3.5776444077
$ws.Range("R2").Value2 = 32.1987996693
$ws.Range("S2").Value2 = 0.09132592689500639
$ws.Range("T2").Value2 = 0.09132592689500639

$ws.Range("G3").Value2 = 1.167885
$ws.Range("H3").Value2 = 3.503655
$ws.Range("I3").Value2 = 0.484739710372708
$ws.Range("J3").Value2 = 0.484739710372708
$ws.Range("O3").Value2 = 0.4156086771445645
$ws.Range("P3").Value2 = 0.4156086771445645
$ws.Range("Q3").Value2 = 7.892167413329999
$ws.Range("R3").Value2 = 71.02950671997
$ws.Range("S3").Value2 = 0.2014620297874405
$ws.Range("T3").Value2 = 0.2014620297874405

$ws.Range("G4").Value2 = 1.167885
$ws.Range("H4").Value2 = 3.503655
$ws.Range("I4").Value2 = 0.484739710372708
$ws.Range("J4").Value2 = 0.484739710372708
$ws.Range("M4").Value2 = 3.493414666666666
$ws.Range("N4").Value2 = 10.480244
$ws.Range("O4").Value2 = 0.214851572590793
$ws.Range("P4").Value2 = 0.214851572590793
$ws.Range("Q4").Value2 = 4.07990658798
$ws.Range("R4").Value2 = 36.71915929182
$ws.Range("S4").Value2 = 0.1041470890707819
$ws.Range("T4").Value2 = 0.1041470890707819

$ws.Range("G5").Value2 = 1.167885
$ws.Range("H5").Value2 = 3.503655
$ws.Range("I5").Value2 = 0.484739710372708
$ws.Range("J5").Value2 = 0.484739710372708
$ws.Range("M5").Value2 = 2.945239333333333
$ws.Range("N5").Value2 = 8.835718
$ws.Range("O5").Value2 = 0.1811377585549322
$ws.Range("P5").Value2 = 0.1811377585549322
$ws.Range("Q5").Value2 = 3.43970083881
$ws.Range("R5").Value2 = 30.95730754929
$ws.Range("S5").Value2 = 0.08780466461947933
$ws.Range("T5").Value2 = 0.08780466461947933

$ws.Range("G6").Value2 = 0.955631
$ws.Range("I6").Value2 = 0.3966420445190934
$ws.Range("J6").Value2 = 0.3966420445190934
$ws.Range("M6").Value2 = 3.063353333333333
$ws.Range("N6").Value2 = 9.190059999999999
$ws.Range("O6").Value2 = 0.1884019917097105
$ws.Range("P6").Value2 = 0.1884019917097105
$ws.Range("Q6").Value2 = 2.927435409286666
$ws.Range("R6").Value2 = 26.34691868358
$ws.Range("S6").Value2 = 0.07472815118320884
$ws.Range("T6").Value2 = 0.07472815118320884

$ws.Range("G7").Value2 = 0.955631
$ws.Range("I7").Value2 = 0.3966420445190934
$ws.Range("J7").Value2 = 0.3966420445190934
$ws.Range("O7").Value2 = 0.4156086771445645
$ws.Range("P7").Value2 = 0.4156086771445645
$ws.Range("Q7").Value2 = 6.457827472198
$ws.Range("R7").Value2 = 58.120447249782
$ws.Range("S7").Value2 = 0.1648478754224958
$ws.Range("T7").Value2 = 0.1648478754224958

$ws.Range("G8").Value2 = 0.955631
$ws.Range("I8").Value2 = 0.3966420445190934
$ws.Range("J8").Value2 = 0.3966420445190934
$ws.Range("M8").Value2 = 3.493414666666666
$ws.Range("N8").Value2 = 10.480244
$ws.Range("O8").Value2 = 0.214851572590793
$ws.Range("P8").Value2 = 0.214851572590793
$ws.Range("Q8").Value2 = 3.338415351321333
$ws.Range("R8").Value2 = 30.045738161892
$ws.Range("S8").Value2 = 0.08521916702055453
$ws.Range("T8").Value2 = 0.08521916702055453

$ws.Range("G9").Value2 = 0.955631
$ws.Range("I9").Value2 = 0.3966420445190934
$ws.Range("J9").Value2 = 0.3966420445190934
$ws.Range("M9").Value2 = 2.945239333333333
$ws.Range("N9").Value2 = 8.835718
$ws.Range("O9").Value2 = 0.1811377585549322
$ws.Range("P9").Value2 = 0.1811377585549322
$ws.Range("Q9").Value2 = 2.814562009352667
$ws.Range("R9").Value2 = 25.331058084174
$ws.Range("S9").Value2 = 0.07184685089283419
$ws.Range("T9").Value2 = 0.07184685089283419

$ws.Range("E10").Value2 = 2
$ws.Range("F10").Value2 = 0.6666666666666666
$ws.Range("G10").Value2 = 0.2857873333333333
$ws.Range("H10").Value2 = 0.857362
$ws.Range("I10").Value2 = 0.1186182451081986
$ws.Range("J10").Value2 = 0.1186182451081986
$ws.Range("M10").Value2 = 3.063353333333333
$ws.Range("N10").Value2 = 9.190059999999999
$ws.Range("O10").Value2 = 0.1884019917097105
$ws.Range("P10").Value2 = 0.1884019917097105
$ws.Range("Q10").Value2 = 0.875467580191111
$ws.Range("R10").Value2 = 7.879208221719999
$ws.Range("S10").Value2 = 0.02234791363149525
$ws.Range("T10").Value2 = 0.02234791363149524

$ws.Range("E11").Value2 = 2
$ws.Range("F11").Value2 = 0.6666666666666666
$ws.Range("G11").Value2 = 0.2857873333333333
$ws.Range("H11").Value2 = 0.857362
$ws.Range("I11").Value2 = 0.1186182451081986
$ws.Range("J11").Value2 = 0.1186182451081986
$ws.Range("O11").Value2 = 0.4156086771445645
$ws.Range("P11").Value2 = 0.4156086771445645
$ws.Range("Q11").Value2 = 1.931253059398667
$ws.Range("R11").Value2 = 17.381277534588
$ws.Range("S11").Value2 = 0.04929877193462814
$ws.Range("T11").Value2 = 0.04929877193462814

$ws.Range("E12").Value2 = 2
$ws.Range("F12").Value2 = 0.6666666666666666
$ws.Range("G12").Value2 = 0.2857873333333333
$ws.Range("H12").Value2 = 0.857362
$ws.Range("I12").Value2 = 0.1186182451081986
$ws.Range("J12").Value2 = 0.1186182451081986
$ws.Range("M12").Value2 = 3.493414666666666
$ws.Range("N12").Value2 = 10.480244
$ws.Range("O12").Value2 = 0.214851572590793
$ws.Range("P12").Value2 = 0.214851572590793
$ws.Range("Q12").Value2 = 0.9983736618142222
$ws.Range("R12").Value2 = 8.985362956327998
$ws.Range("S12").Value2 = 0.02548531649945662
$ws.Range("T12").Value2 = 0.02548531649945661

$ws.Range("E13").Value2 = 2
$ws.Range("F13").Value2 = 0.6666666666666666
$ws.Range("G13").Value2 = 0.2857873333333333
$ws.Range("H13").Value2 = 0.857362
$ws.Range("I13").Value2 = 0.1186182451081986
$ws.Range("J13").Value2 = 0.1186182451081986
$ws.Range("M13").Value2 = 2.945239333333333
$ws.Range("N13").Value2 = 8.835718
$ws.Range("O13").Value2 = 0.1811377585549322
$ws.Range("P13").Value2 = 0.1811377585549322
$ws.Range("Q13").Value2 = 0.8417120951017778
$ws.Range("R13").Value2 = 7.575408855916
$ws.Range("S13").Value2 = 0.02148624304261865
$ws.Range("T13").Value2 = 0.02148624304261865
